$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.043.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "'1.668.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'216.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").Value = "'0.5100"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.2658"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "'0.06402"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").Value = "'21.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Value = "'0.07443"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "'1.671.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "'4.506"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "'0.5847"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "'0.000008569"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "'64.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "'26.081.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "'4.944"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("D21").Value = "'193.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.70%  "
$ws.Range("D22").Value = "'6.215"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'144.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'7.615"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("D27").Value = "'15.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value = "'0.06468"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.71%  "
$ws.Range("D29").Value = "'1.334"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").Value = "'3.551"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").Value = "'3.519"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'1.652"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'0.6105"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").Value = "'2.368"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'6.255"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.01603"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").Value = "'1.094.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").Value = "'0.8609"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "'100.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.22%  "
$ws.Range("D44").Value = "'1.818.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").Value = "'0.00000000115"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'56.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "'1.011"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "'0.05240"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'0.4284"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").Value = "'6.050"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.36%  "
